$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 19500
$ws.Range("I21").Value = 19500
$ws.Range("K21").Value = 19500
$ws.Range("M21").Value = -19032

$ws.Range("H23").Value = 19500
$ws.Range("I23").Value = 19500
$ws.Range("K23").Value = 19500
$ws.Range("M23").Value = -19266

$ws.Range("H113").Value = 20000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 20000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 20000
$ws.Range("M113").Value = $null
$ws.Range("N113").Value = -26508

$ws.Range("H129").Value = 400
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5551.3335
$ws.Range("I45").Value = 5551.3335
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 5551.3335
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -5174.3335
$ws.Range("N45").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 1249.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 880
$ws.Range("I134").Value = 876.6667
$ws.Range("K134").Value = 2630.0001
$ws.Range("M134").Value = -95.0001000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3930.5557
$ws.Range("I5").Value = 3483.2856
$ws.Range("J5").Value = 5496
$ws.Range("K5").Value = 10449.8568
$ws.Range("L5").Value = 16488
$ws.Range("M5").Value = -10337.8568
$ws.Range("N5").Value = -16712

$ws.Range("H15").Value = 113.55556
$ws.Range("J15").Value = 146
$ws.Range("L15").Value = 438
$ws.Range("N15").Value = -718

$ws.Range("H26").Value = 616.75
$ws.Range("J26").Value = 2242
$ws.Range("L26").Value = 6726
$ws.Range("N26").Value = -7302

$ws.Range("H44").Value = 919.5
$ws.Range("I44").Value = 216.66667
$ws.Range("J44").Value = 1003.84
$ws.Range("K44").Value = 650.00001
$ws.Range("L44").Value = 3011.52
$ws.Range("M44").Value = -252.00001
$ws.Range("N44").Value = -3807.52

$ws.Range("H46").Value = 5498.5
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 5498.5
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 16495.5
$ws.Range("M46").Value = $null
$ws.Range("N46").Value = -16677.5

$ws.Range("H49").Value = 4000
$ws.Range("J49").Value = 4000
$ws.Range("L49").Value = 12000
$ws.Range("N49").Value = -12312

$ws.Range("H86").Value = 461.66666
$ws.Range("I86").Value = 592.5
$ws.Range("J86").Value = 200
$ws.Range("K86").Value = 1777.5
$ws.Range("L86").Value = 600
$ws.Range("M86").Value = -591.5
$ws.Range("N86").Value = -2972

$ws.Range("H89").Value = 461.66666
$ws.Range("I89").Value = 592.5
$ws.Range("J89").Value = 200
$ws.Range("K89").Value = 5332.5
$ws.Range("L89").Value = 1800
$ws.Range("M89").Value = 595.5
$ws.Range("N89").Value = -13656

$ws.Range("H92").Value = 550
$ws.Range("I92").Value = 550
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1650
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -402
$ws.Range("N92").Value = $null

$ws.Range("H135").Value = 3930.5557
$ws.Range("I135").Value = 3483.2856
$ws.Range("J135").Value = 5496
$ws.Range("K135").Value = 31349.5704
$ws.Range("L135").Value = 49464
$ws.Range("M135").Value = -28814.5704
$ws.Range("N135").Value = -54534

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 77.42308
$ws.Range("I2").Value = 58.61111
$ws.Range("J2").Value = 119.75
$ws.Range("K2").Value = 58.61111
$ws.Range("L2").Value = 119.75
$ws.Range("M2").Value = 54.38889
$ws.Range("N2").Value = -345.75

$ws.Range("H24").Value = 19999
$ws.Range("J24").Value = 19999
$ws.Range("L24").Value = 19999
$ws.Range("N24").Value = -20345

$ws.Range("H107").Value = 1879.2
$ws.Range("J107").Value = 1879.2
$ws.Range("L107").Value = 1879.2
$ws.Range("N107").Value = -5719.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").Value = $null

$ws.Range("H32").Value = 5193.3335
$ws.Range("I32").Value = 1232
$ws.Range("K32").Value = 1232
$ws.Range("M32").Value = -915

$ws.Range("H42").Value = 32999.332
$ws.Range("I42").Value = 32498
$ws.Range("J42").Value = 33250
$ws.Range("K42").Value = 32498
$ws.Range("L42").Value = 33250
$ws.Range("M42").Value = -31935
$ws.Range("N42").Value = -34376

$ws.Range("H49").Value = 32999.332
$ws.Range("I49").Value = 32498
$ws.Range("J49").Value = 33250
$ws.Range("K49").Value = 32498
$ws.Range("L49").Value = 33250
$ws.Range("M49").Value = -32351
$ws.Range("N49").Value = -33544

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 2500
$ws.Range("J3").Value = 2500
$ws.Range("L3").Value = 2500
$ws.Range("N3").Value = -2728

$ws.Range("H33").Value = 16000
$ws.Range("I33").Value = 16000
$ws.Range("K33").Value = 16000
$ws.Range("M33").Value = -15750

$ws.Range("H36").Value = 16000
$ws.Range("I36").Value = 16000
$ws.Range("K36").Value = 16000
$ws.Range("M36").Value = -15750

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").Value = $null

$ws.Range("H81").Value = 2622
$ws.Range("I81").Value = 2622
$ws.Range("K81").Value = 5244
$ws.Range("M81").Value = -4183

$ws.Range("H82").Value = 25000
$ws.Range("I82").Value = 25000
$ws.Range("K82").Value = 25000
$ws.Range("M82").Value = -24617

$ws.Range("H84").Value = 2622
$ws.Range("I84").Value = 2622
$ws.Range("K84").Value = 26220
$ws.Range("M84").Value = -20916

$ws.Range("H85").Value = 25000
$ws.Range("I85").Value = 25000
$ws.Range("K85").Value = 25000
$ws.Range("M85").Value = -23674

$ws.Range("H107").Value = 1074.8889
$ws.Range("I107").Value = 812.4286
$ws.Range("K107").Value = 2437.2858
$ws.Range("M107").Value = -517.2857999999997

$ws.Range("H132").Value = 1290.0834
$ws.Range("I132").Value = 1099
$ws.Range("J132").Value = 1672.25
$ws.Range("K132").Value = 3297
$ws.Range("L132").Value = 5016.75
$ws.Range("M132").Value = -767
$ws.Range("N132").Value = -10076.75
